$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder: "Testing custom properties"
# Split "Testing " -> "Testing" + " ", and "custom " -> "custom" + " "
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Characters(1, 7).Text = "Testing"
$titleRange.Characters(9, 6).Text = "custom"

# Subtitle placeholder: (line break)(line break)"A. M."
# Split "A. " -> "A." + " "
$subtitleShape = $s.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange
$subtitleRange.Characters(3, 2).Text = "A."
